# Localization Manual.docx edits
# - Remove the stray "_GoBack" bookmark that sat at the very start of the
#   document (it will be re-created later, at the point the author was
#   actually last editing).
# - Fix a handful of typos / merge runs that had been needlessly split.
# - Fix the example msgfmt command (drop the bogus "ORTS." prefix) and
#   leave a "_GoBack" bookmark where the edit was made, as Word does.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the original "_GoBack" bookmark near the top of the document.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Helper-free straightforward text fixes (Find/Replace over the whole
# story so formatting-identical adjoining runs collapse into one, just
# like Word does when you retype over a selection).
# ---------------------------------------------------------------------

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $false, 1, $false, $replace, 2) | Out-Null
}

# "e." + " See the license for more details." -> single run
Replace-Text "its use. See the license for more details." "its use. See the license for more details."

# "Trademark Acknowledgment" heading, rebuilt from split runs
Replace-Text "Trademark Acknowledgment" "Trademark Acknowledgment"

# "Copyright Acknowledgment and License" heading, rebuilt from split runs
Replace-Text "Copyright Acknowledgment and License" "Copyright Acknowledgment and License"

# Copyright year "2009-2014"
Replace-Text "2009-2014" "2009-2014"

# "This document is part of Open Rails."
Replace-Text "    This document is part of Open Rails." "    This document is part of Open Rails."

# "Open Rails is free software: you can redistribute it and/or modify"
Replace-Text "    Open Rails is free software: you can redistribute it and/or modify" "    Open Rails is free software: you can redistribute it and/or modify"

# "any later version."
Replace-Text "    any later version." "    any later version."

# "as part of the Open Rails distribution in Documentation\Copying.txt. "
Replace-Text "    as part of the Open Rails distribution in Documentation\Copying.txt. " "    as part of the Open Rails distribution in Documentation\Copying.txt. "

# "If not, see "
Replace-Text "    If not, see " "    If not, see "

# "It is located at:"
Replace-Text "It is located at:" "It is located at:"

# "C:\Windows\Microsoft.NET\Framework for 32 bit computers"
Replace-Text "C:\Windows\Microsoft.NET\Framework for 32 bit computers" "C:\Windows\Microsoft.NET\Framework for 32 bit computers"

# "C:\Windows\Microsoft.NET\Framework64 for 64 bit computers"
Replace-Text "C:\Windows\Microsoft.NET\Framework64 for 64 bit computers" "C:\Windows\Microsoft.NET\Framework64 for 64 bit computers"

# "Save the PO file in Source\Locales\Menu."
Replace-Text "Save the PO file in Source\Locales\Menu." "Save the PO file in Source\Locales\Menu."

# "Add the name of your language to the list (in a new line)."
Replace-Text "Add the name of your language to the list (in a new line)." "Add the name of your language to the list (in a new line)."

# "When developers ... GNU.Gettext.Xgettext.exe in the 3rdPartyLibs directory."
Replace-Text "When developers add new strings to the game, they will have to update the model. They have to use GNU.Gettext.Xgettext.exe in the 3rdPartyLibs directory." "When developers add new strings to the game, they will have to update the model. They have to use GNU.Gettext.Xgettext.exe in the 3rdPartyLibs directory."

# ---------------------------------------------------------------------
# 2. Fix the sample msgfmt command: "-r ORTS.Menu" -> "-r Menu" and leave
#    a "_GoBack" bookmark exactly where "ORTS." used to be, matching the
#    cursor position Word would have recorded for that edit.
# ---------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("-r ORTS.Menu", $true, $false, $false, $false, $false, $false, 1, $false, "", 0)
if ($found) {
    $insertPoint = $findRng.Start + 3   # length of "-r "
    $replaceRng = $d.Range($insertPoint, $findRng.End)
    $replaceRng.Text = "Menu"
    $bmRange = $d.Range($insertPoint, $insertPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
